# Update NATMI TPM-derived ligand/receptor edge weights for Sema6b-Plxna2.
#
# The upstream scripts were re-run with new TPM input, which changed the
# per-cluster average expression of the ligand (Sema6b) and the receptor
# (Plxna2). Everything else on the sheet (average/total expression,
# derived-specificity, and edge-weight columns) is recomputed from those
# two per-cluster tables, exactly like the original NATMI pipeline does:
#   total   = average * expressing-cell-count
#   spec    = own value / sum(value over all clusters)
#   edge avg/total/spec = ligand value * matching receptor value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusterOrder = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Neutrophils", "Resolving-Mac")

# Refreshed ligand (Sema6b) average expression per sending cluster.
$ligandAvg = @{
    "ECs"               = 14.69933033333333
    "FAPs"              = 3.463629333333333
    "Inflammatory-Mac"  = 3.397405666666666
    "MuSCs"             = 0.442328
    "Neutrophils"       = 11.94113466666666
    "Resolving-Mac"     = 1.969077666666667
}

# Refreshed receptor (Plxna2) average expression per target cluster.
$receptorAvg = @{
    "ECs"               = 39.62362533333334
    "FAPs"              = 2.247453666666667
    "Inflammatory-Mac"  = 0.2964306666666667
    "MuSCs"             = 2.447182
    "Neutrophils"       = 0.4200656666666667
    "Resolving-Mac"     = 0.9146693333333333
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- figure out which clusters' ligand/receptor averages actually moved,
#     using column G/M of the first row that references each cluster, so
#     that columns whose inputs are unchanged are left completely alone
#     (avoids introducing floating point noise on untouched cells). ---
$ligandChanged = @{}
$receptorChanged = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target = $ws.Cells.Item($r, 4).Value2

    if (-not $ligandChanged.ContainsKey($sending)) {
        $oldG = $ws.Cells.Item($r, 7).Value2
        $ligandChanged[$sending] = [Math]::Abs($oldG - $ligandAvg[$sending]) -gt 0.000000000001
    }
    if (-not $receptorChanged.ContainsKey($target)) {
        $oldM = $ws.Cells.Item($r, 13).Value2
        $receptorChanged[$target] = [Math]::Abs($oldM - $receptorAvg[$target]) -gt 0.000000000001
    }
}

$cellCount = 3

$ligandTotal = @{}
$receptorTotal = @{}
foreach ($k in $clusterOrder) {
    $ligandTotal[$k] = $ligandAvg[$k] * $cellCount
    $receptorTotal[$k] = $receptorAvg[$k] * $cellCount
}

$sumLigandAvg = 0.0
$sumLigandTotal = 0.0
$sumReceptorAvg = 0.0
$sumReceptorTotal = 0.0
foreach ($k in $clusterOrder) {
    $sumLigandAvg = $sumLigandAvg + $ligandAvg[$k]
    $sumLigandTotal = $sumLigandTotal + $ligandTotal[$k]
    $sumReceptorAvg = $sumReceptorAvg + $receptorAvg[$k]
    $sumReceptorTotal = $sumReceptorTotal + $receptorTotal[$k]
}

$ligandSpecAvg = @{}
$ligandSpecTotal = @{}
$receptorSpecAvg = @{}
$receptorSpecTotal = @{}
foreach ($k in $clusterOrder) {
    $ligandSpecAvg[$k] = $ligandAvg[$k] / $sumLigandAvg
    $ligandSpecTotal[$k] = $ligandTotal[$k] / $sumLigandTotal
    $receptorSpecAvg[$k] = $receptorAvg[$k] / $sumReceptorAvg
    $receptorSpecTotal[$k] = $receptorTotal[$k] / $sumReceptorTotal
}

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target = $ws.Cells.Item($r, 4).Value2

    $i = $ligandSpecAvg[$sending]
    $j = $ligandSpecTotal[$sending]
    $o = $receptorSpecAvg[$target]
    $p = $receptorSpecTotal[$target]

    if ($ligandChanged[$sending]) {
        $ws.Cells.Item($r, 7).Value = $ligandAvg[$sending]
        $ws.Cells.Item($r, 8).Value = $ligandTotal[$sending]
    }
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j

    if ($receptorChanged[$target]) {
        $ws.Cells.Item($r, 13).Value = $receptorAvg[$target]
        $ws.Cells.Item($r, 14).Value = $receptorTotal[$target]
    }
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p

    if ($ligandChanged[$sending] -or $receptorChanged[$target]) {
        $ws.Cells.Item($r, 17).Value = $ligandAvg[$sending] * $receptorAvg[$target]
        $ws.Cells.Item($r, 18).Value = $ligandTotal[$sending] * $receptorTotal[$target]
    }
    $ws.Cells.Item($r, 19).Value = $i * $o
    $ws.Cells.Item($r, 20).Value = $j * $p
}
